$d = $word.ActiveDocument

# The first paragraph ends with "corregida ". We want to insert a brand new
# paragraph containing "Mi segunda línea de código " right after it, before
# the existing (empty) second paragraph.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()

# The newly created paragraph is now Paragraphs(2); fill in its text.
$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "Mi segunda línea de código "
